$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q0)
$ws.Range("B2").Value = 0.3349553228780944
$ws.Range("C2").Value = 0.3367524377135537
$ws.Range("D2").Value = 0.1614313307467032
$ws.Range("E2").Value = 0.4017851798495101
$ws.Range("F2").Value = 0.2302686848532658

# Row 3 (Q1)
$ws.Range("B3").Value = 0.1821895727347061
$ws.Range("C3").Value = 0.1995134267841846
$ws.Range("D3").Value = 0.06522829055202861
$ws.Range("E3").Value = 0.2553982978643918
$ws.Range("F3").Value = 0.1886656364482292

# Row 4 (Q2)
$ws.Range("B4").Value = 0.2077349899393594
$ws.Range("C4").Value = 0.2383528991286697
$ws.Range("D4").Value = 0.09280927092432688
$ws.Range("E4").Value = 0.3046461405045645
$ws.Range("F4").Value = 0.2441035310172005

# Row 5 (Q3)
$ws.Range("B5").Value = 0.1283418663009961
$ws.Range("C5").Value = 0.1283418663009961
$ws.Range("D5").Value = 0.02714726664902342
$ws.Range("E5").Value = 0.1647642760097692
$ws.Range("F5").Value = 0.146120717240237
